# Weekly refresh of fruit/vegetable prices: Fruta, Vega Monumental Concepción - Frambuesa
# Rows 3-19 hold weekly price records; this update shifts the date/quality/
# price figures for most weeks (rows 12 & 13 are unchanged).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3
$ws.Range("D3").Value = 44216
$ws.Range("N3").Value = 3500
$ws.Range("O3").Value = 4000
$ws.Range("P3").Value = 3750
$ws.Range("S3").Value = 1875

# Row 4
$ws.Range("D4").Value = 44216
$ws.Range("L4").Value = "Segunda"
$ws.Range("M4").Value = 100
$ws.Range("O4").Value = 3000
$ws.Range("P4").Value = 3000
$ws.Range("S4").Value = 1500

# Row 5
$ws.Range("D5").Value = 44574
$ws.Range("L5").Value = "Primera"
$ws.Range("M5").Value = 200
$ws.Range("N5").Value = 7000
$ws.Range("O5").Value = 8000
$ws.Range("P5").Value = 7500
$ws.Range("S5").Value = 3750

# Row 6
$ws.Range("D6").Value = 44574
$ws.Range("L6").Value = "Segunda"
$ws.Range("N6").Value = 6000
$ws.Range("O6").Value = 6000
$ws.Range("P6").Value = 6000
$ws.Range("S6").Value = 3000

# Row 7
$ws.Range("L7").Value = "Primera"
$ws.Range("N7").Value = 10000
$ws.Range("O7").Value = 10000
$ws.Range("P7").Value = 10000
$ws.Range("S7").Value = 5000

# Row 8
$ws.Range("D8").Value = 44532
$ws.Range("L8").Value = "Segunda"
$ws.Range("N8").Value = 8000
$ws.Range("P8").Value = 8000
$ws.Range("S8").Value = 4000

# Row 9
$ws.Range("D9").Value = 44195
$ws.Range("N9").Value = 3000
$ws.Range("O9").Value = 3500
$ws.Range("P9").Value = 3250
$ws.Range("S9").Value = 1625

# Row 10
$ws.Range("D10").Value = 44195
$ws.Range("N10").Value = 2500
$ws.Range("O10").Value = 2500
$ws.Range("P10").Value = 2500
$ws.Range("S10").Value = 1250

# Row 11
$ws.Range("D11").Value = 44617
$ws.Range("N11").Value = 6000
$ws.Range("O11").Value = 7000
$ws.Range("P11").Value = 6500
$ws.Range("S11").Value = 3250

# Rows 12 & 13: unchanged

# Row 14
$ws.Range("D14").Value = 44917
$ws.Range("N14").Value = 7000
$ws.Range("O14").Value = 7500
$ws.Range("P14").Value = 7250
$ws.Range("S14").Value = 3625

# Row 15
$ws.Range("D15").Value = 44602
$ws.Range("L15").Value = "Primera"
$ws.Range("M15").Value = 200
$ws.Range("N15").Value = 6000
$ws.Range("O15").Value = 7000
$ws.Range("P15").Value = 6500
$ws.Range("S15").Value = 3250

# Row 16
$ws.Range("L16").Value = "Segunda"
$ws.Range("M16").Value = 100
$ws.Range("N16").Value = 5000
$ws.Range("O16").Value = 5000
$ws.Range("P16").Value = 5000
$ws.Range("S16").Value = 2500

# Row 17
$ws.Range("D17").Value = 44559
$ws.Range("L17").Value = "Primera"
$ws.Range("M17").Value = 200
$ws.Range("N17").Value = 6000
$ws.Range("O17").Value = 7000
$ws.Range("P17").Value = 6500
$ws.Range("S17").Value = 3250

# Row 18
$ws.Range("D18").Value = 44559
$ws.Range("L18").Value = "Segunda"
$ws.Range("M18").Value = 100
$ws.Range("N18").Value = 5000
$ws.Range("O18").Value = 5000
$ws.Range("P18").Value = 5000
$ws.Range("S18").Value = 2500

# Row 19
$ws.Range("D19").Value = 44944
$ws.Range("L19").Value = "Primera"
$ws.Range("N19").Value = 7000
$ws.Range("O19").Value = 8000
$ws.Range("P19").Value = 7500
$ws.Range("S19").Value = 3750
